# Apply the "enh: update plies written" change to the evolution workbook.
# Sheet "results": row 14 (previously STDEV.S of Q6:Q10 .. Y6:Y10) is repurposed
# to hold the AVERAGE of the same ranges and labeled "Average n of plies (i1-i4)".
# A new row 15 stores the raw "N of plies 5th iteration" values (copied from row 10),
# and a new row 16 restores the original STDEV calculation (now labeled "std" again),
# using STDEV.P for column Q and STDEV.S for columns R:Y, matching the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# --- Row 14: switch from STDEV.S(...) to AVERAGE(...) and relabel ---
$ws.Range("P14").Value = "Average n of plies (i1-i4)"

$ws.Range("Q14").Formula = "=AVERAGE(Q6:Q10)"
$ws.Range("R14").Formula = "=AVERAGE(R6:R10)"
$ws.Range("S14").Formula = "=AVERAGE(S6:S10)"
$ws.Range("T14").Formula = "=AVERAGE(T6:T10)"
$ws.Range("U14").Formula = "=AVERAGE(U6:U10)"
$ws.Range("V14").Formula = "=AVERAGE(V6:V10)"
$ws.Range("W14").Formula = "=AVERAGE(W6:W10)"
$ws.Range("X14").Formula = "=AVERAGE(X6:X10)"
$ws.Range("Y14").Formula = "=AVERAGE(Y6:Y10)"

# --- Row 15 (new): literal "N of plies 5th iteration" values (copy of row 10) ---
$ws.Range("P15").Value = "N of plies 5th iteration"
$ws.Range("Q15").Value = 156
$ws.Range("R15").Value = 119
$ws.Range("S15").Value = 79
$ws.Range("T15").Value = 82
$ws.Range("U15").Value = 126
$ws.Range("V15").Value = 163
$ws.Range("W15").Value = 164
$ws.Range("X15").Value = 124
$ws.Range("Y15").Value = 87

# --- Row 16 (new): the STDEV formulas that used to live in row 14 ---
$ws.Range("Q16:Y16").NumberFormat = "0.0"
$ws.Range("P16").Value = "std"
$ws.Range("Q16").Formula = "=_xlfn.STDEV.P(Q6:Q10)"
$ws.Range("R16").Formula = "=_xlfn.STDEV.S(R6:R10)"
$ws.Range("S16").Formula = "=_xlfn.STDEV.S(S6:S10)"
$ws.Range("T16").Formula = "=_xlfn.STDEV.S(T6:T10)"
$ws.Range("U16").Formula = "=_xlfn.STDEV.S(U6:U10)"
$ws.Range("V16").Formula = "=_xlfn.STDEV.S(V6:V10)"
$ws.Range("W16").Formula = "=_xlfn.STDEV.S(W6:W10)"
$ws.Range("X16").Formula = "=_xlfn.STDEV.S(X6:X10)"
$ws.Range("Y16").Formula = "=_xlfn.STDEV.S(Y6:Y10)"

# --- Column P width: widen column 16 (P) and drop the old "best fit" sizing ---
$ws.Columns.Item(16).ColumnWidth = 20.45

# --- View state: scroll the frozen pane back to the top and move the selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("L66").Select()
